# Atualizacao de bases das ligas, do dia: 17-02-2024 as 22:47
#
# A handful of match rows in the "Kazakhstan Premier League" sheet had been
# written against the wrong fixture id. This re-aligns each row's betting
# data (HomeTeam/AwayTeam/odds/etc., columns B and F:AC) with its correct
# match while leaving the row's own sequence number (column A, "id") and the
# other row-identity columns (Div, Div Original Name) untouched - i.e. it
# swaps/rotates the payload between the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the per-match payload (everything except A/C/D/E, which
# never change for these rows).
$cols = @('B','F','G','H','I','J','K','L','M','N','O','P','Q','R','S','T','U','V','W','X','Y','Z','AA','AB','AC')

function Get-RowData([int]$row, $cols) {
    $data = @{}
    foreach ($col in $cols) {
        $data[$col] = $ws.Range("$col$row").Value2
    }
    return $data
}

function Set-RowData([int]$row, $data) {
    foreach ($col in $data.Keys) {
        $ws.Range("$col$row").Value = $data[$col]
    }
}

function Swap-Rows([int]$rowA, [int]$rowB, $cols) {
    $dataA = Get-RowData $rowA $cols
    $dataB = Get-RowData $rowB $cols
    Set-RowData $rowA $dataB
    Set-RowData $rowB $dataA
}

# Simple pairwise swaps.
Swap-Rows 16 17 $cols
Swap-Rows 19 20 $cols
Swap-Rows 27 28 $cols
Swap-Rows 38 39 $cols
Swap-Rows 88 89 $cols

# Rows 177/178/179 form a 3-way rotation rather than a plain swap:
#   after(177) = before(179)
#   after(178) = before(177)
#   after(179) = before(178)
$d177 = Get-RowData 177 $cols
$d178 = Get-RowData 178 $cols
$d179 = Get-RowData 179 $cols

Set-RowData 177 $d179
Set-RowData 178 $d177
Set-RowData 179 $d178
